# Rename the first 8 worksheet tabs to fix the "plutonium" (钚) naming:
#  - "回收钸" (recycle-Cu-like typo) -> "钚回收" (plutonium recovery)
#  - "铀钸镄" -> "铀钚镄" (fix mis-typed character 钸 -> 钚 in Uranium/Plutonium/Fermium)

$wb = $excel.ActiveWorkbook

$renames = @(
    @{ Index = 1; Name = "最大化发电-无废料-钚回收-允许转化-有红石" },
    @{ Index = 2; Name = "最大化发电-无废料-钚回收-允许转化-无红石" },
    @{ Index = 3; Name = "最大化发电-无废料-钚回收-无转化-有红石" },
    @{ Index = 4; Name = "最大化发电-无废料-钚回收-无转化-无红石" },
    @{ Index = 5; Name = "最大化发电-无废料-铀钚镄-允许转化-有红石" },
    @{ Index = 6; Name = "最大化发电-无废料-铀钚镄-允许转化-无红石" },
    @{ Index = 7; Name = "最大化发电-无废料-铀钚镄-无转化-有红石" },
    @{ Index = 8; Name = "最大化发电-无废料-铀钚镄-无转化-无红石" }
)

foreach ($r in $renames) {
    $ws = $wb.Worksheets.Item($r.Index)
    $ws.Name = $r.Name
}
